$wb = $excel.ActiveWorkbook

# --- STEP 1: Worker Host Configuration (Cluster Configuration sheet) ---
$wsCluster = $wb.Worksheets.Item("Cluster Configuration")

# RAM (GB) per worker host: 256 -> 128
$wsCluster.Range("D8").Value = 128

# CPU (cores) per worker host: 48 -> 20
$wsCluster.Range("D9").Value = 20

# HDD (TB) per worker host: 36 -> 12
$wsCluster.Range("D10").Value = 12

# Physical Cores to Vcores Multiplier: 4 -> 2
$wsCluster.Range("E29").Value = 2

# --- Update the active cell / selection on each sheet to reflect where the
#     author was last working, then leave "Cluster Configuration" as the
#     selected/active tab. ---
$wsYarn = $wb.Worksheets.Item("YARN Configuration")
$wsYarn.Activate()
$wsYarn.Range("E46").Select()

$wsMapReduce = $wb.Worksheets.Item("MapReduce Configuration")
$wsMapReduce.Activate()
$wsMapReduce.Range("C47").Select()

$wsCluster.Activate()
$wsCluster.Range("I46").Select()
